# This script applies targeted cell-value updates to Sheet1 of the workbook,
# matching the changes described in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: cell address, new value
$updates = @{
    "J18"  = '2.7/8"'
    "J19"  = '2.7/8"'
    "J34"  = '2.7/8"'
    "K35"  = 'V-630'
    "J36"  = '3.1/2"'
    "J61"  = '2000 PSI'
    "L68"  = 'NETZSCH'
    "J81"  = '06'
    "J82"  = '06'
    "L84"  = 'NETZSCH'
    "L85"  = 'NETZSCH'
    "J96"  = '7/8"'
    "J103" = 'P8'
    "K110" = 'EL'
    "J112" = '7/8"'
    "J113" = '7/8"'
    "J114" = '7/8"'
    "J116" = '7/8"'
    "J132" = '27'
    "J137" = '70'
    "J152" = '20'
    "K153" = 'EL'
    "J154" = '20'
    "J155" = '15'
    "J161" = '3.1/2"'
    "L161" = 'API 5B'
    "J162" = '3.1/2"'
    "L162" = 'API 5B'
    "J163" = '90'
    "L163" = 'NTZ'
    "J166" = '2.7/8"'
    "L166" = 'API 5B'
    "J177" = '10'
    "J181" = '2.7/8"'
    "L181" = 'API 5B'
    "L182" = 'TI'
    "J184" = '5.1/2"'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so purely-numeric-looking strings (e.g. "06", "40")
    # keep their original text representation instead of becoming numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

$wb.Save()
